# "browser closed, added dates"
# Fill in the start/end period dates for the two companies on the
# "Параметры" sheet (columns C = Начало периода, D = Окончание периода).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ООО "СЕГМЕНТЭНЕРГО": set the format + value first...
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = "01/23/2001"
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Value = "12/24/2018"

# ...then reuse the exact same formatting for row 3 - ООО "ЕВРОТРАНС" by
# copying the format down, so both rows share one style record.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("C3").Value = "01/12/1990"
$ws.Range("D3").Value = "09/23/2018"

$excel.CutCopyMode = 0

# Move the active selection to D4, matching where the user left off editing.
[void]$ws.Range("D4").Select()
